# scopes.xlsx update — add a "printers" column (F/G) to the scope table,
# nudge the active selection, and resize the header row / new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing "x" marker cell that was blank before.
$ws.Range("C4").Value = "x"

# New "x" marker + note cell, styled to match the other centered marker
# cells in this row (D4/E4 use centered horizontal+vertical alignment).
$markers = $ws.Range("F4:G4")
$markers.HorizontalAlignment = -4108   # xlCenter
$markers.VerticalAlignment = -4108     # xlCenter
$ws.Range("F4").Value = "x"
$ws.Range("G4").Value = "substituir com o laravel scout"

# New header cell for the "printers" scope group.
$ws.Range("F1").Value = "printers"

# Column sizing for the newly-used F/G columns.
$ws.Columns.Item(6).ColumnWidth = 2.6
$ws.Columns.Item(7).ColumnWidth = 32.1

# Header row height shrank slightly to fit the new layout.
$ws.Rows.Item(1).RowHeight = 65.25

# Selection moved from H6 to H9.
$null = $ws.Range("H9").Select()
